$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "US290"
$ws.Range("C3").Value = "IH610"

$ws.Range("C4").Select()
